$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.273.46"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Value = "3.696.05"
$ws.Range("E3").Value = "  -3.63%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.41"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.88"
$ws.Range("E6").Value = "  +11.10%  "
$ws.Range("D7").Value = "3.687.00"
$ws.Range("E7").Value = "  -3.83%  "
$ws.Range("E8").Value = "  -5.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.723"
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("E11").Value = "  -5.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.67"
$ws.Range("E12").Value = "  +6.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000291"
$ws.Range("E13").Value = "  -8.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.43"
$ws.Range("E14").Value = "  -5.72%  "
$ws.Range("D15").Value = "4.281.12"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("D16").Value = "3.697.43"
$ws.Range("E16").Value = "  -4.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.40"
$ws.Range("E17").Value = "  -5.52%  "
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("E19").Value = "  -5.62%  "
$ws.Range("E20").Value = "  -6.59%  "
$ws.Range("D21").Value = "68.019.35"
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "410.59"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.47"
$ws.Range("E23").Value = "  -4.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.58"
$ws.Range("E24").Value = "  -5.28%  "
$ws.Range("E25").Value = "  -6.62%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.81"
$ws.Range("E26").Value = "  -6.28%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.84"
$ws.Range("E28").Value = "  -6.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.05"
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.48"
$ws.Range("E30").Value = "  -5.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.71"
$ws.Range("E31").Value = "  -5.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.43"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.53"
$ws.Range("E33").Value = "  -6.50%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "44.21"
$ws.Range("E34").Value = "  -10.89%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.117"
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.87"
$ws.Range("E36").Value = "  -5.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "590.94"
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("D38").Value = "0.0₃0890"
$ws.Range("E38").Value = "  -7.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.402"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.73"
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0436"
$ws.Range("E45").Value = "  -5.88%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  -9.60%  "
$ws.Range("E47").Value = "  -7.34%  "
$ws.Range("D48").Value = "2.785.93"
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("E49").Value = "  -5.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.67"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.12"
$ws.Range("E51").Value = "  -7.43%  "
